$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values (initial/final state ranges) and column B values (type of process),
# per row, expressed as formulas that evaluate to string literals (mirrors the
# authoring tool re-writing shared-string cells as ="..." text formulas).
$colA = @(
    "initial_final_states",
    "1-2",
    "2-3",
    "3-4",
    "4-5",
    "5-6",
    "6-7",
    "7-8",
    "8-9",
    "9-10",
    "10-11",
    "11-12",
    "12-13"
)

$colB = @(
    "type_of_process",
    "isobaric",
    "isochoric",
    "isothermal",
    "isothermal",
    "adiabatic",
    "isochoric",
    "adiabatic",
    "isobaric",
    "isobaric",
    "isothermal",
    "adiabatic",
    "isochoric"
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Formula = '="' + $colA[$i] + '"'
    $ws.Cells.Item($row, 2).Formula = '="' + $colB[$i] + '"'
}

# Update the selected cell shown in the sheet view.
$ws.Range("B3").Select()

# Mark the two error-checker hints the author dismissed by hand in Excel
# ("Ignore Error" on the little cell-corner warning): an inconsistent-formula
# warning on B7 and a two-digit-text-year warning on A13.
try { $ws.Range("B7").Errors.Item(4).Value = $true } catch {}
try { $ws.Range("A13").Errors.Item(9).Value = $true } catch {}

$wb.Save()
